$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New parts added to the inventory list (rows 37 and 38).
# Shared-string order in the target file: "Futaba 3003 Servo" is referenced
# first (row 38 / col A) so it gets the lowest new index, then
# "Tower Pro SG90" (row 37 / col A) and " 9 grams" (row 37 / col C).
$ws.Range("A38").Value = "Futaba 3003 Servo"
$ws.Range("A37").Value = "Tower Pro SG90"
$ws.Range("C37").Value = " 9 grams"

$ws.Range("B37").Value = 8
$ws.Range("B38").Value = 10

# Copy the formatting (cell styles / borders / fill) of the last existing
# data row down onto the two new rows, matching the rest of the table.
$ws.Range("A36:E36").Copy()
$ws.Range("A37:E37").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A36:E36").Copy()
$ws.Range("A38:E38").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# D37/E37, C38/D38/E38 stay blank (formatted only, no value) just like D36/E36.
$ws.Range("D37").ClearContents()
$ws.Range("E37").ClearContents()
$ws.Range("C38").ClearContents()
$ws.Range("D38").ClearContents()
$ws.Range("E38").ClearContents()

# Match the saved view: scrolled so row 10 is at the top, with F38 selected.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F38").Select()
